$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $ws.Range($cell).NumberFormat = "@"
    $ws.Range($cell).Value = $val
    $ws.Range($cell).Style = "Normal"
}

$ws.Range("D2").Value = "28.886.58"
$ws.Range("E2").Value = "  -1.86%  "
$ws.Range("D3").Value = "1.825.29"
$ws.Range("E3").Value = "  -2.25%  "
Set-TextValue "D4" "0.9989"
$ws.Range("E4").Value = "  -0.22%  "
Set-TextValue "D5" "239.12"
$ws.Range("E5").Value = "  -1.75%  "
Set-TextValue "D6" "0.6902"
$ws.Range("E6").Value = "  -1.99%  "
Set-TextValue "D7" "0.9995"
$ws.Range("E7").Value = "  -0.21%  "
Set-TextValue "D8" "0.07614"
$ws.Range("E8").Value = "  -2.98%  "
Set-TextValue "D9" "0.3014"
$ws.Range("E9").Value = "  -3.97%  "
Set-TextValue "D10" "23.38"
$ws.Range("E10").Value = "  -4.39%  "
Set-TextValue "D11" "0.07720"
$ws.Range("E11").Value = "  -3.82%  "
$ws.Range("D12").Value = "1.817.74"
$ws.Range("E12").Value = "  -3.97%  "
Set-TextValue "D13" "5.040"
$ws.Range("E13").Value = "  -2.93%  "
Set-TextValue "D14" "89.94"
$ws.Range("E14").Value = "  -3.68%  "
Set-TextValue "D15" "0.6712"
$ws.Range("E15").Value = "  -4.15%  "
Set-TextValue "D16" "6.403"
$ws.Range("E16").Value = "  -0.69%  "
Set-TextValue "D17" "0.000008269"
$ws.Range("E17").Value = "  -0.58%  "
$ws.Range("D18").Value = "28.870.66"
$ws.Range("E18").Value = "  -2.14%  "
Set-TextValue "D19" "242.63"
$ws.Range("E19").Value = "  -4.55%  "
$ws.Range("D20").Value = "2.080.14"
$ws.Range("E20").Value = "  -2.86%  "
Set-TextValue "D21" "12.61"
$ws.Range("E21").Value = "  -3.90%  "
Set-TextValue "D22" "0.9995"
Set-TextValue "D23" "7.372"
$ws.Range("E23").Value = "  -2.92%  "
Set-TextValue "D24" "0.9990"
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("E25").Value = "  -5.25%  "
Set-TextValue "D26" "160.22"
Set-TextValue "D27" "8.702"
$ws.Range("E27").Value = "  -3.52%  "
Set-TextValue "D28" "18.14"
$ws.Range("E28").Value = "  -3.29%  "
Set-TextValue "D29" "1.528"
$ws.Range("E29").Value = "  +1.81%  "
Set-TextValue "D30" "4.184"
$ws.Range("E30").Value = "  -3.20%  "
Set-TextValue "D31" "4.130"
$ws.Range("E31").Value = "  -2.96%  "
$ws.Range("E32").Value = "  -0.87%  "
Set-TextValue "D33" "0.05096"
$ws.Range("E33").Value = "  -3.83%  "
Set-TextValue "D34" "0.7484"
$ws.Range("E34").Value = "  +0.37%  "
Set-TextValue "D35" "1.810"
$ws.Range("E35").Value = "  -3.98%  "
Set-TextValue "D36" "1.140"
$ws.Range("E36").Value = "  -2.03%  "
Set-TextValue "D37" "2.683"
$ws.Range("E37").Value = "  -1.31%  "
Set-TextValue "D38" "0.01830"
$ws.Range("E38").Value = "  -2.09%  "
$ws.Range("D39").Value = "1.196.83"
$ws.Range("E39").Value = "  -4.93%  "
Set-TextValue "D40" "2.674"
$ws.Range("E40").Value = "  -2.42%  "
Set-TextValue "D41" "0.9137"
$ws.Range("E41").Value = "  +1.73%  "
Set-TextValue "D42" "108.31"
$ws.Range("E42").Value = "  -0.19%  "
Set-TextValue "D43" "0.9988"
$ws.Range("E43").Value = "  -0.25%  "
$ws.Range("D44").Value = "1.979.93"
$ws.Range("E44").Value = "  -2.92%  "
$ws.Range("E45").Value = "  -5.55%  "
Set-TextValue "D46" "0.5149"
$ws.Range("E46").Value = "  -0.79%  "
Set-TextValue "D47" "9.441"
$ws.Range("E47").Value = "  -0.56%  "
Set-TextValue "D48" "5.230"
$ws.Range("E48").Value = "  -12.19%  "
Set-TextValue "D49" "1.725"
$ws.Range("E49").Value = "  -3.69%  "
Set-TextValue "D50" "62.18"
$ws.Range("E50").Value = "  -12.53%  "
Set-TextValue "D51" "0.4190"
$ws.Range("E51").Value = "  -2.60%  "
